$wb = $excel.ActiveWorkbook

$wsCreate  = $wb.Worksheets.Item("Create")
$wsEdit    = $wb.Worksheets.Item("Edit")
$wsDelete  = $wb.Worksheets.Item("Delete")

# --- Content edits -------------------------------------------------

# Create sheet: update color code, add a new start-time row
$wsCreate.Range("C2").Value = "#0e101e"
$wsCreate.Range("A3").Value = "'08:00:00"
$wsCreate.Range("A3").NumberFormat = "h:mm:ss"

# Edit sheet: update color codes
$wsEdit.Range("D2").Value = "#0e101e"
$wsEdit.Range("E2").Value = "#1e0e16"

# Delete sheet: update start time + color codes
$wsDelete.Range("A2").Value = "08:48:00"
$wsDelete.Range("C2").Value = "#0e101e"
$wsDelete.Range("D2").Value = "#1e0e16"

# --- Selection / active sheet --------------------------------------

[void]$wsDelete.Range("A2").Select()
[void]$wsEdit.Range("A2").Select()
[void]$wsCreate.Activate()
[void]$wsCreate.Range("A3").Select()
